# "Ajout de légendes et d'une option pour quitter"
# Adds two new journal entries (rows 28-29) to the work-log table and
# extends the Tableau1 table with 11 additional blank rows (30-40) so the
# table keeps its usual pool of ready-to-fill rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Before touching anything, clone the still-blank template look of
#    row 28 onto the 11 rows that will become the new blank tail of the
#    table (rows 30-40), so they keep the unfilled-row formatting.
# ---------------------------------------------------------------------
$ws.Range("E28:M28").Copy() | Out-Null
$ws.Range("E30:M40").PasteSpecial(-4122) | Out-Null

for ($r = 30; $r -le 40; $r++) {
  $ws.Range("H$r").Formula = '=IF(ISBLANK(Tableau1[[#This Row],[Heure Début]]),"",Tableau1[[#This Row],[Heure fin]]-Tableau1[[#This Row],[Heure Début]])'
}

# ---------------------------------------------------------------------
# 2) Turn the previously-blank template row 28 into a real entry, and
#    build a fresh entry in row 29, both copying the formatting already
#    used by the existing data rows (keeps date/time/duration/text
#    cell styles identical to the rest of the journal).
# ---------------------------------------------------------------------
$ws.Range("E27:M27").Copy()
$ws.Range("E28:M29").PasteSpecial(-4122)

$ws.Range("E28").Value = 44266
$ws.Range("F28").Value = 0.74305555555555547
$ws.Range("G28").Value = 0.75694444444444453
$ws.Range("H28").Formula = '=IF(ISBLANK(Tableau1[[#This Row],[Heure Début]]),"",Tableau1[[#This Row],[Heure fin]]-Tableau1[[#This Row],[Heure Début]])'
$ws.Range("I28").Value = "Développement"
$ws.Range("J28").Value = "Légendes"
$ws.Range("K28").Value = "Domicile"
$ws.Range("L28").Value = "Ajout de légende sur les côté de la grille"

$ws.Range("E29").Value = 44266
$ws.Range("F29").Value = 0.76388888888888884
$ws.Range("G29").Value = 0.77777777777777779
$ws.Range("H29").Formula = '=IF(ISBLANK(Tableau1[[#This Row],[Heure Début]]),"",Tableau1[[#This Row],[Heure fin]]-Tableau1[[#This Row],[Heure Début]])'
$ws.Range("I29").Value = "Développement"
$ws.Range("L29").Value = "Ajout d'une option pour quitter pendant le jeux"
$ws.Range("J29").Value = "ajout d'une option de quitter"
$ws.Range("K29").Value = "Domicile"

# ---------------------------------------------------------------------
# 3) Resize the Tableau1 list-object + its AutoFilter to cover the newly
#    added rows, E5:M28 -> E5:M40.
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("E5:M40"))

# ---------------------------------------------------------------------
# 4) Leave the selection on the next empty row, as it was when the
#    author saved the file.
# ---------------------------------------------------------------------
$ws.Range("J30").Select()
